$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Step 1: Move rows 4-13 of sheet1 ("工作表1") to the end of sheet2 ("工作表2") as rows 4-13,
#     adding two new columns F/G with new traffic-count data.

$ws2.Cells.Item(4, 1).Value = "大度路出城"
$ws2.Cells.Item(4, 2).Value = 3967
$ws2.Cells.Item(4, 3).Value = 4196
$ws2.Cells.Item(4, 4).Value = 8322
$ws2.Cells.Item(4, 5).Value = 15698
$ws2.Cells.Item(4, 6).Value = 10465
$ws2.Cells.Item(4, 7).Value = 10225

$ws2.Cells.Item(5, 1).Value = "大度路進城"
$ws2.Cells.Item(5, 2).Value = 5018
$ws2.Cells.Item(5, 3).Value = 3682
$ws2.Cells.Item(5, 4).Value = 15697
$ws2.Cells.Item(5, 5).Value = 8122
$ws2.Cells.Item(5, 6).Value = 10226
$ws2.Cells.Item(5, 7).Value = 8126

$ws2.Cells.Item(6, 1).Value = "重陽橋進城"
$ws2.Cells.Item(6, 2).Value = 5480
$ws2.Cells.Item(6, 3).Value = 2758
$ws2.Cells.Item(6, 4).Value = 7536
$ws2.Cells.Item(6, 5).Value = 5268
$ws2.Cells.Item(6, 6).Value = 7534
$ws2.Cells.Item(6, 7).Value = 6668

$ws2.Cells.Item(7, 1).Value = "重陽橋出城"
$ws2.Cells.Item(7, 2).Value = 2370
$ws2.Cells.Item(7, 3).Value = 3975
$ws2.Cells.Item(7, 4).Value = 5268
$ws2.Cells.Item(7, 5).Value = 7536
$ws2.Cells.Item(7, 6).Value = 6669
$ws2.Cells.Item(7, 7).Value = 5831

$ws2.Cells.Item(8, 1).Value = "臺北橋進城"
$ws2.Cells.Item(8, 2).Value = 7719
$ws2.Cells.Item(8, 3).Value = 3181
$ws2.Cells.Item(8, 4).Value = 7607
$ws2.Cells.Item(8, 5).Value = 7609
$ws2.Cells.Item(8, 6).Value = 8111
$ws2.Cells.Item(8, 7).Value = 6542

$ws2.Cells.Item(9, 1).Value = "臺北橋出城"
$ws2.Cells.Item(9, 2).Value = 2041
$ws2.Cells.Item(9, 3).Value = 6321
$ws2.Cells.Item(9, 4).Value = 7609
$ws2.Cells.Item(9, 5).Value = 7607
$ws2.Cells.Item(9, 6).Value = 6544
$ws2.Cells.Item(9, 7).Value = 5202

$ws2.Cells.Item(10, 1).Value = "忠孝橋進城"
$ws2.Cells.Item(10, 2).Value = 7570
$ws2.Cells.Item(10, 3).Value = 4516
$ws2.Cells.Item(10, 4).Value = 5797
$ws2.Cells.Item(10, 5).Value = 5295
$ws2.Cells.Item(10, 6).Value = 5799
$ws2.Cells.Item(10, 7).Value = 6511

$ws2.Cells.Item(11, 1).Value = "忠孝橋出城"
$ws2.Cells.Item(11, 2).Value = 3521
$ws2.Cells.Item(11, 3).Value = 6780
$ws2.Cells.Item(11, 4).Value = 5295
$ws2.Cells.Item(11, 5).Value = 5797
$ws2.Cells.Item(11, 6).Value = 6511
$ws2.Cells.Item(11, 7).Value = 5798

$ws2.Cells.Item(12, 1).Value = "中興橋進城"
$ws2.Cells.Item(12, 2).Value = 4696
$ws2.Cells.Item(12, 3).Value = 3026
$ws2.Cells.Item(12, 4).Value = 5783
$ws2.Cells.Item(12, 5).Value = 5772
$ws2.Cells.Item(12, 6).Value = 16895
$ws2.Cells.Item(12, 7).Value = 8155

$ws2.Cells.Item(13, 1).Value = "中興橋出城"
$ws2.Cells.Item(13, 2).Value = 2183
$ws2.Cells.Item(13, 3).Value = 3274
$ws2.Cells.Item(13, 4).Value = 5772
$ws2.Cells.Item(13, 5).Value = 5783
$ws2.Cells.Item(13, 6).Value = 8166
$ws2.Cells.Item(13, 7).Value = 7603

# --- Step 2: Remove old rows 4-13 from sheet1 (data now lives in sheet2); remaining rows shift up.
$ws1.Rows("4:13").Delete()

# --- Step 3: Append 12 new rows (new survey links) to the bottom of sheet1, rows 22-33.
$ws1.Cells.Item(22, 1).Value = "民族(-中山北)東向"
$ws1.Cells.Item(22, 2).Value = 4278
$ws1.Cells.Item(22, 3).Value = 1743
$ws1.Cells.Item(22, 4).Value = 6627
$ws1.Cells.Item(22, 5).Value = 8346

$ws1.Cells.Item(23, 1).Value = "民族(中山北-)西向"
$ws1.Cells.Item(23, 2).Value = 1122
$ws1.Cells.Item(23, 3).Value = 2408
$ws1.Cells.Item(23, 4).Value = 8346
$ws1.Cells.Item(23, 5).Value = 6627

$ws1.Cells.Item(24, 1).Value = "民權(-中山北)東向"
$ws1.Cells.Item(24, 2).Value = 5089
$ws1.Cells.Item(24, 3).Value = 2695
$ws1.Cells.Item(24, 4).Value = 15501
$ws1.Cells.Item(24, 5).Value = 8353

$ws1.Cells.Item(25, 1).Value = "民權(中山北-)西向"
$ws1.Cells.Item(25, 2).Value = 2029
$ws1.Cells.Item(25, 3).Value = 4467
$ws1.Cells.Item(25, 4).Value = 8347
$ws1.Cells.Item(25, 5).Value = 15500

$ws1.Cells.Item(26, 1).Value = "民生(-中山北)東向"
$ws1.Cells.Item(26, 2).Value = 2204
$ws1.Cells.Item(26, 3).Value = 1761
$ws1.Cells.Item(26, 4).Value = 15499
$ws1.Cells.Item(26, 5).Value = 8360

$ws1.Cells.Item(27, 1).Value = "民生(中山北-)西向"
$ws1.Cells.Item(27, 2).Value = 947
$ws1.Cells.Item(27, 3).Value = 1340
$ws1.Cells.Item(27, 4).Value = 8359
$ws1.Cells.Item(27, 5).Value = 15498

$ws1.Cells.Item(28, 1).Value = "南京(-中山北)東向"
$ws1.Cells.Item(28, 2).Value = 2452
$ws1.Cells.Item(28, 3).Value = 1368
$ws1.Cells.Item(28, 4).Value = 9898
$ws1.Cells.Item(28, 5).Value = 8365

$ws1.Cells.Item(29, 1).Value = "南京(中山北-)西向"
$ws1.Cells.Item(29, 2).Value = 1034
$ws1.Cells.Item(29, 3).Value = 2480
$ws1.Cells.Item(29, 4).Value = 8365
$ws1.Cells.Item(29, 5).Value = 9898

$ws1.Cells.Item(30, 1).Value = "市民(-中山北)東向"
$ws1.Cells.Item(30, 2).Value = 3234
$ws1.Cells.Item(30, 3).Value = 1873
$ws1.Cells.Item(30, 4).Value = 15442
$ws1.Cells.Item(30, 5).Value = 8376

$ws1.Cells.Item(31, 1).Value = "市民(中山北-)西向"
$ws1.Cells.Item(31, 2).Value = 1317
$ws1.Cells.Item(31, 3).Value = 2726
$ws1.Cells.Item(31, 4).Value = 8371
$ws1.Cells.Item(31, 5).Value = 15443

$ws1.Cells.Item(32, 1).Value = "忠孝(-中山北)東向"
$ws1.Cells.Item(32, 2).Value = 1909
$ws1.Cells.Item(32, 3).Value = 1181
$ws1.Cells.Item(32, 4).Value = 8382
$ws1.Cells.Item(32, 5).Value = 7928

$ws1.Cells.Item(33, 1).Value = "忠孝(中山北-)西向"
$ws1.Cells.Item(33, 2).Value = 1472
$ws1.Cells.Item(33, 3).Value = 1916
$ws1.Cells.Item(33, 4).Value = 7928
$ws1.Cells.Item(33, 5).Value = 8383

# --- Step 4: Restore sheet view selections to match the saved state.
$ws2.Activate()
$ws2.Range("A14:E15").Select()
$ws1.Activate()
$ws1.Range("D30:E30").Select()
try {
    $excel.ActiveWindow.ScrollRow = 28
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}

Write-Output "Edit complete"
